$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep text formatting so values like
# "0.2765" or "1.793.38" are not auto-converted to numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.633.05'
$ws.Range("E2").Value = '  -2.18%  '

$ws.Range("D3").Value = '1.793.38'
$ws.Range("E3").Value = '  -2.04%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '231.75'
$ws.Range("E5").Value = '  -1.75%  '

$ws.Range("E6").Value = '  -2.51%  '

$ws.Range("E7").Value = '  -0.15%  '

$ws.Range("D8").Value = '0.2765'
$ws.Range("E8").Value = '  -1.00%  '

$ws.Range("E9").Value = '  -4.25%  '

$ws.Range("E10").Value = '  -1.58%  '

$ws.Range("D11").Value = '0.07526'
$ws.Range("E11").Value = '  -1.65%  '

$ws.Range("D12").Value = '1.794.82'
$ws.Range("E12").Value = '  -2.02%  '

$ws.Range("D13").Value = '4.782'
$ws.Range("E13").Value = '  -0.16%  '

$ws.Range("D14").Value = '0.6126'
$ws.Range("E14").Value = '  -2.44%  '

$ws.Range("D15").Value = '2.035.43'
$ws.Range("E15").Value = '  -2.05%  '

$ws.Range("E16").Value = '  -4.66%  '

$ws.Range("D17").Value = '0.000008837'
$ws.Range("E17").Value = '  -9.37%  '

$ws.Range("D18").Value = '28.611.06'
$ws.Range("E18").Value = '  -2.16%  '

$ws.Range("D19").Value = '5.410'
$ws.Range("E19").Value = '  -7.10%  '

$ws.Range("E20").Value = '  -0.06%  '

$ws.Range("D21").Value = '208.66'
$ws.Range("E21").Value = '  -6.89%  '

$ws.Range("E22").Value = '  -1.94%  '

$ws.Range("D23").Value = '6.832'
$ws.Range("E23").Value = '  -2.44%  '

$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("D25").Value = '152.91'
$ws.Range("E25").Value = '  -2.30%  '

$ws.Range("D26").Value = '8.112'
$ws.Range("E26").Value = '  +1.60%  '

$ws.Range("D27").Value = '0.1259'
$ws.Range("E27").Value = '  -3.20%  '

$ws.Range("D28").Value = '16.37'
$ws.Range("E28").Value = '  -1.45%  '

$ws.Range("D29").Value = '1.415'
$ws.Range("E29").Value = '  -3.70%  '

$ws.Range("D30").Value = '0.06237'
$ws.Range("E30").Value = '  -6.38%  '

$ws.Range("E31").Value = '  -1.80%  '

$ws.Range("D32").Value = '3.805'
$ws.Range("E32").Value = '  +0.39%  '

$ws.Range("D33").Value = '3.780'
$ws.Range("E33").Value = '  -1.67%  '

$ws.Range("E34").Value = '  +0.83%  '

$ws.Range("D35").Value = '1.047'
$ws.Range("E35").Value = '  -5.35%  '

$ws.Range("D36").Value = '0.6391'
$ws.Range("E36").Value = '  -1.03%  '

$ws.Range("D37").Value = '2.499'
$ws.Range("E37").Value = '  -2.00%  '

$ws.Range("D38").Value = '2.712'
$ws.Range("E38").Value = '  -0.90%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01694'
$ws.Range("E39").Value = '  -3.70%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '6.360'
$ws.Range("E40").Value = '  -2.55%  '

$ws.Range("D41").Value = '1.141.62'
$ws.Range("E41").Value = '  -5.91%  '

$ws.Range("D42").Value = '0.8748'
$ws.Range("E42").Value = '  -2.81%  '

$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("D44").Value = '100.23'
$ws.Range("E44").Value = '  -0.19%  '

$ws.Range("D45").Value = '1.945.80'
$ws.Range("E45").Value = '  -2.20%  '

$ws.Range("D46").Value = '59.90'
$ws.Range("E46").Value = '  -4.28%  '

$ws.Range("D47").Value = '0.00000000110'
$ws.Range("E47").Value = '  -4.23%  '

$ws.Range("E48").Value = '  +0.38%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.360'
$ws.Range("E49").Value = '  -1.76%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.05458'
$ws.Range("E50").Value = '  -0.81%  '

$ws.Range("E51").Value = '  -1.60%  '

# Restore default cell style (no explicit style index) while keeping text values.
$ws.Range("D2:D51").Style = "Normal"
